$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:B1").ClearContents()
$ws.Rows("1").AutoFit()
$ws.Columns("A:B").Select()
